$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0004542274691630155

$ws.Range("A3").Value = 0.0001223599974764511
$ws.Range("C3").Value = 43.07692337036133
$ws.Range("D3").Value = 21.538461685180664

$ws.Range("A4").Value = 0.00010876038868445903
$ws.Range("C4").Value = 47.30769348144531
$ws.Range("D4").Value = 23.653844833374023

$ws.Range("A5").Value = 0.00006291375029832125
$ws.Range("C5").Value = 47.5
$ws.Range("D5").Value = 41.98219299316406

$ws.Range("A6").Value = 0.000038863268855493516
$ws.Range("C6").Value = 47.30769348144531
$ws.Range("D6").Value = 23.653844833374023

$ws.Range("A7").Value = 0.00003810634734691121
$ws.Range("C7").Value = 47.30769348144531
$ws.Range("D7").Value = 23.653844833374023

$ws.Range("A8").Value = 0.000030950770451454446
$ws.Range("C8").Value = 43.07692337036133
$ws.Range("D8").Value = 21.538461685180664

$ws.Range("A9").Value = 0.000027236250389250927
$ws.Range("C9").Value = 51.730770111083984
$ws.Range("D9").Value = 25.94953727722168

$ws.Range("A10").Value = 0.00001579326817591209
$ws.Range("C10").Value = 43.269229888916016
$ws.Range("D10").Value = 21.698057174682617

$ws.Range("A11").Value = 0.000006102692623244366
$ws.Range("C11").Value = 47.30769348144531
$ws.Range("D11").Value = 23.376941680908203

$ws.Range("D12").Value = 2.274712085723877
